$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Noelia changed her Codewars username: noeliacarrion -> nihal.99
$ws.Range("A13").Value = "nihal.99"

# 2. Add a new kata column D: "filter-out-the-geese"
# Copy the formatting of the neighboring header cell (C1) onto the new
# header cell (D1) before writing its text.
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("D1").Value = "filter-out-the-geese"

# Completion status (TRUE/FALSE) for each user (rows 2-21) for the new kata
$completed = @($true,$true,$true,$true,$true,$false,$false,$true,$true,$true,$true,$true,$false,$true,$true,$true,$true,$true,$true,$true)

for ($i = 0; $i -lt $completed.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $completed[$i]
}
